# [EGSVC-36] Modified create employee feature
# Adds two new sample rows (employee3/employee4, assignment3/assignment4,
# JurisdictionList3/JurisdictionList4) to the eisTestData workbook used by
# the create-employee functional test.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: employeeDetails
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("employeeDetails")

# Row 4 - employee3
$ws1.Range("A4").Value = "employee3"
$ws1.Range("B4").Value = "Deputation"
$ws1.Range("C4").Value = "DECEASED"
$ws1.Range("D4").NumberFormat = "DD/MM/YY"
$ws1.Range("D4").Value = [DateTime]"1988-01-01"
$ws1.Range("E4").Value = "Male"
$ws1.Range("F4").Value = "MARRIED"
$ws1.Range("G4").Value = "Yes"
$ws1.Range("H4").NumberFormat = "@"
$ws1.Range("H4").Value = "7777777777"
$ws1.Range("I4").Value = "Municipal Office Rd, N.R.Peta, Near Appollo Hospital, Kurnool, Andhra Pradesh "
$ws1.Range("J4").Value = "Kurnool"
$ws1.Range("L4").NumberFormat = "@"
$ws1.Range("L4").Value = "02/01/2013"

# Row 5 - employee4
$ws1.Range("A5").Value = "employee4"
$ws1.Range("B5").Value = "Outsourced"
$ws1.Range("C5").Value = "EMPLOYED"
$ws1.Range("D5").NumberFormat = "DD/MM/YY"
$ws1.Range("D5").Value = [DateTime]"1991-01-01"
$ws1.Range("E5").Value = "Male"
$ws1.Range("F5").Value = "UNMARRIED"
$ws1.Range("G5").Value = "Yes"
$ws1.Range("H5").NumberFormat = "@"
$ws1.Range("H5").Value = "6666666666"
$ws1.Range("I5").Value = "Municipal Office Rd, N.R.Peta, Near Appollo Hospital, Kurnool, Andhra Pradesh "
$ws1.Range("J5").Value = "Kurnool"
$ws1.Range("L5").NumberFormat = "@"
$ws1.Range("L5").Value = "02/01/2014"

# ---------------------------------------------------------------------
# Sheet: assignmentDetails
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("assignmentDetails")

# Row 4 - assignment3
$ws2.Range("A4").Value = "assignment3"
$ws2.Range("B4").Value = "No"
$ws2.Range("C4").NumberFormat = "DD/MM/YY"
$ws2.Range("C4").Value = [DateTime]"2017-05-05"
$ws2.Range("D4").NumberFormat = "DD/MM/YY"
$ws2.Range("D4").Value = [DateTime]"2017-05-05"
$ws2.Range("E4").Value = "ACCOUNTS"
$ws2.Range("F4").Value = "Accounts Officer"
$ws2.Range("G4").Value = "ACC_ACC_1"

# Row 5 - assignment4
$ws2.Range("A5").Value = "assignment4"
$ws2.Range("B5").Value = "Yes"
$ws2.Range("C5").NumberFormat = "DD/MM/YY"
$ws2.Range("C5").Value = [DateTime]"2017-06-05"
$ws2.Range("D5").NumberFormat = "DD/MM/YY"
$ws2.Range("D5").Value = [DateTime]"2017-06-05"
$ws2.Range("F5").Value = "Accounts Officer"
$ws2.Range("G5").Value = "ACC_Accounts_Officer_001"

# ---------------------------------------------------------------------
# Sheet: jurisdictionList
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("jurisdictionList")

# Row 4 - JurisdictionList3
$ws3.Range("A4").Value = "JurisdictionList3"
$ws3.Range("B4").Value = "Ward"
$ws3.Range("C4").Value = "Election Ward No. 2"

# Row 5 - JurisdictionList4
$ws3.Range("A5").Value = "JurisdictionList4"
$ws3.Range("B5").Value = "Ward"
$ws3.Range("C5").Value = "Election Ward No. 3"
